$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-27 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-28 Monday", 2) | Out-Null
$d.Content.Find.Execute("75×28=2100", $true, $false, $false, $false, $false, $true, 1, $false, "87×25=2175", 2) | Out-Null
$d.Content.Find.Execute("53×16=848", $true, $false, $false, $false, $false, $true, 1, $false, "48×68=3264", 2) | Out-Null
$d.Content.Find.Execute("41×71=2911", $true, $false, $false, $false, $false, $true, 1, $false, "79×60=4740", 2) | Out-Null
$d.Content.Find.Execute("93×87=8091", $true, $false, $false, $false, $false, $true, 1, $false, "85×21=1785", 2) | Out-Null
$d.Content.Find.Execute("55×57=3135", $true, $false, $false, $false, $false, $true, 1, $false, "80×34=2720", 2) | Out-Null
$d.Content.Find.Execute("91×68=6188", $true, $false, $false, $false, $false, $true, 1, $false, "29×35=1015", 2) | Out-Null
$d.Content.Find.Execute("20×83=1660", $true, $false, $false, $false, $false, $true, 1, $false, "88×43=3784", 2) | Out-Null
$d.Content.Find.Execute("57×96=5472", $true, $false, $false, $false, $false, $true, 1, $false, "94×37=3478", 2) | Out-Null
$d.Content.Find.Execute("97×16=1552", $true, $false, $false, $false, $false, $true, 1, $false, "76×92=6992", 2) | Out-Null
$d.Content.Find.Execute("35×86=3010", $true, $false, $false, $false, $false, $true, 1, $false, "19×83=1577", 2) | Out-Null
$d.Content.Find.Execute("89×89=7921", $true, $false, $false, $false, $false, $true, 1, $false, "74×73=5402", 2) | Out-Null
$d.Content.Find.Execute("89×38=3382", $true, $false, $false, $false, $false, $true, 1, $false, "34×82=2788", 2) | Out-Null
$d.Content.Find.Execute("99×76=7524", $true, $false, $false, $false, $false, $true, 1, $false, "28×67=1876", 2) | Out-Null
$d.Content.Find.Execute("35×55=1925", $true, $false, $false, $false, $false, $true, 1, $false, "39×90=3510", 2) | Out-Null
$d.Content.Find.Execute("15×94=1410", $true, $false, $false, $false, $false, $true, 1, $false, "86×56=4816", 2) | Out-Null
$d.Content.Find.Execute("91×99=9009", $true, $false, $false, $false, $false, $true, 1, $false, "70×14=980", 2) | Out-Null
$d.Content.Find.Execute("87×98=8526", $true, $false, $false, $false, $false, $true, 1, $false, "34×32=1088", 2) | Out-Null
$d.Content.Find.Execute("24×23=552", $true, $false, $false, $false, $false, $true, 1, $false, "84×53=4452", 2) | Out-Null
$d.Content.Find.Execute("49×66=3234", $true, $false, $false, $false, $false, $true, 1, $false, "59×24=1416", 2) | Out-Null
$d.Content.Find.Execute("59×60=3540", $true, $false, $false, $false, $false, $true, 1, $false, "15×35=525", 2) | Out-Null
$d.Content.Find.Execute("49×52=2548", $true, $false, $false, $false, $false, $true, 1, $false, "36×43=1548", 2) | Out-Null
$d.Content.Find.Execute("27×24=648", $true, $false, $false, $false, $false, $true, 1, $false, "21×75=1575", 2) | Out-Null
$d.Content.Find.Execute("67×34=2278", $true, $false, $false, $false, $false, $true, 1, $false, "29×11=319", 2) | Out-Null
$d.Content.Find.Execute("85×35=2975", $true, $false, $false, $false, $false, $true, 1, $false, "44×91=4004", 2) | Out-Null
$d.Content.Find.Execute("75×24=1800", $true, $false, $false, $false, $false, $true, 1, $false, "27×65=1755", 2) | Out-Null
